$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 419, shifting existing rows 419-511 down to 420-512.
$ws.Rows.Item(419).Insert()

# Populate the newly inserted row with the latest weekly price-report entry.
$ws.Cells.Item(419, 1).Value = 6
$ws.Cells.Item(419, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(419, 3).Value = "Metropolitana"
$ws.Cells.Item(419, 4).Value = 44641
$ws.Cells.Item(419, 5).Value = 13
$ws.Cells.Item(419, 6).Value = 100112012
$ws.Cells.Item(419, 7).Value = "Espinaca"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 510
$ws.Cells.Item(419, 11).Value = 5500
$ws.Cells.Item(419, 12).Value = 6000
$ws.Cells.Item(419, 13).Value = 5725
$ws.Cells.Item(419, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(419, 15).Value = "Región Metropolitana"
$ws.Cells.Item(419, 16).Value = 572
$ws.Cells.Item(419, 17).Value = 10
$ws.Cells.Item(419, 18).Value = "Hortaliza"
